$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1852.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1852.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 1852.5
$ws.Range("N62").Value = -3100.5
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 1852.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1852.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 9262.5
$ws.Range("N65").Value = -15502.5
$ws.Range("M65").ClearContents()
$ws.Range("H113").Value = 3650
$ws.Range("I113").Value = 3650
$ws.Range("K113").Value = 3650
$ws.Range("M113").Value = -396
$ws.Range("H137").Value = 35724036
$ws.Range("I137").Value = 55558016
$ws.Range("J137").Value = 22870
$ws.Range("K137").Value = 166674048
$ws.Range("L137").Value = 68610
$ws.Range("M137").Value = -166671498
$ws.Range("N137").Value = -73710

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 657.05
$ws.Range("I32").Value = 657.05
$ws.Range("K32").Value = 657.05
$ws.Range("M32").Value = -370.05
$ws.Range("H45").Value = 2443.65
$ws.Range("I45").Value = 1871.6
$ws.Range("K45").Value = 1871.6
$ws.Range("M45").Value = -1494.6
$ws.Range("H74").Value = 1164407.6
$ws.Range("I74").Value = 1987928.2
$ws.Range("J74").Value = 11478.6
$ws.Range("K74").Value = 1987928.2
$ws.Range("L74").Value = 11478.6
$ws.Range("M74").Value = -1987054.2
$ws.Range("N74").Value = -13226.6
$ws.Range("H77").Value = 1164407.6
$ws.Range("I77").Value = 1987928.2
$ws.Range("J77").Value = 11478.6
$ws.Range("K77").Value = 9939641
$ws.Range("L77").Value = 57393
$ws.Range("M77").Value = -9935273
$ws.Range("N77").Value = -66129
$ws.Range("H132").Value = 3996.875
$ws.Range("I132").Value = 3997.5
$ws.Range("J132").Value = 3995
$ws.Range("K132").Value = 11992.5
$ws.Range("L132").Value = 11985
$ws.Range("M132").Value = -9462.5
$ws.Range("N132").Value = -17045
$ws.Range("H135").Value = 80000
$ws.Range("J135").Value = 80000
$ws.Range("L135").Value = 80000
$ws.Range("N135").Value = -90140
$ws.Range("H139").Value = 91510.73
$ws.Range("J139").Value = 91510.73
$ws.Range("L139").Value = 91510.73
$ws.Range("N139").Value = -101790.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 64799.8
$ws.Range("J81").Value = 64799.8
$ws.Range("L81").Value = 64799.8
$ws.Range("N81").Value = -66921.8
$ws.Range("H84").Value = 64799.8
$ws.Range("J84").Value = 64799.8
$ws.Range("L84").Value = 194399.4
$ws.Range("N84").Value = -205007.4
$ws.Range("H134").Value = 6675595
$ws.Range("I134").Value = 5056.8096
$ws.Range("K134").Value = 15170.4288
$ws.Range("M134").Value = -12635.4288
$ws.Range("H135").Value = 60929.6
$ws.Range("J135").Value = 60929.6
$ws.Range("L135").Value = 60929.6
$ws.Range("N135").Value = -71069.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1896581.1
$ws.Range("I31").Value = 2086019.2
$ws.Range("J31").Value = 2200
$ws.Range("K31").Value = 2086019.2
$ws.Range("L31").Value = 2200
$ws.Range("M31").Value = -2085724.2
$ws.Range("N31").Value = -2790
$ws.Range("H34").Value = 1896581.1
$ws.Range("I34").Value = 2086019.2
$ws.Range("J34").Value = 2200
$ws.Range("K34").Value = 2086019.2
$ws.Range("L34").Value = 2200
$ws.Range("M34").Value = -2085817.2
$ws.Range("N34").Value = -2604
$ws.Range("H86").Value = 17453.04
$ws.Range("J86").Value = 18976.1
$ws.Range("L86").Value = 18976.1
$ws.Range("N86").Value = -21222.1
$ws.Range("H89").Value = 17453.04
$ws.Range("J89").Value = 18976.1
$ws.Range("L89").Value = 94880.5
$ws.Range("N89").Value = -106112.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3250.8
$ws.Range("J3").Value = 5500
$ws.Range("L3").Value = 16500
$ws.Range("N3").Value = -16724
$ws.Range("H9").Value = 94903.12
$ws.Range("I9").Value = 201570.2
$ws.Range("J9").Value = 50458.5
$ws.Range("K9").Value = 604710.6000000001
$ws.Range("L9").Value = 151375.5
$ws.Range("M9").Value = -604486.6000000001
$ws.Range("N9").Value = -151823.5
$ws.Range("H17").Value = 389.25
$ws.Range("I17").Value = 513
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 1539
$ws.Range("L17").Value = 54
$ws.Range("M17").Value = -1370
$ws.Range("N17").Value = -392
$ws.Range("H63").Value = 512
$ws.Range("I63").Value = 512
$ws.Range("K63").Value = 1536
$ws.Range("M63").Value = -787
$ws.Range("H66").Value = 512
$ws.Range("I66").Value = 512
$ws.Range("K66").Value = 4608
$ws.Range("M66").Value = -864
$ws.Range("H131").Value = 1605.0526
$ws.Range("J131").Value = 1605.0526
$ws.Range("L131").Value = 4815.1578
$ws.Range("N131").Value = -14895.1578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 773.1111
$ws.Range("I107").Value = 1034.6666
$ws.Range("J107").Value = 250
$ws.Range("K107").Value = 1034.6666
$ws.Range("L107").Value = 250
$ws.Range("M107").Value = 885.3334
$ws.Range("N107").Value = -4090
$ws.Range("H122").Value = 3650.4285
$ws.Range("I122").Value = 2749.8333
$ws.Range("J122").Value = 4325.875
$ws.Range("K122").Value = 8249.499899999999
$ws.Range("L122").Value = 12977.625
$ws.Range("M122").Value = -5799.499899999999
$ws.Range("N122").Value = -17877.625
$ws.Range("H132").Value = 43229.9
$ws.Range("I132").Value = 46922.11
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 140766.33
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -138236.33
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2499
$ws.Range("I40").Value = 2373.75
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2373.75
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2237.75
$ws.Range("N40").Value = -3272
$ws.Range("H46").Value = 3872.8333
$ws.Range("J46").Value = 4286.875
$ws.Range("L46").Value = 4286.875
$ws.Range("N46").Value = -4662.875
$ws.Range("H61").Value = 11743.333
$ws.Range("I61").Value = 11901.818
$ws.Range("K61").Value = 11901.818
$ws.Range("M61").Value = -11699.818
$ws.Range("H74").Value = 47179.727
$ws.Range("I74").Value = 37056.4
$ws.Range("J74").Value = 55615.832
$ws.Range("K74").Value = 37056.4
$ws.Range("L74").Value = 55615.832
$ws.Range("M74").Value = -36058.4
$ws.Range("N74").Value = -57611.832
$ws.Range("H77").Value = 47179.727
$ws.Range("I77").Value = 37056.4
$ws.Range("J77").Value = 55615.832
$ws.Range("K77").Value = 111169.2
$ws.Range("L77").Value = 166847.496
$ws.Range("M77").Value = -106177.2
$ws.Range("N77").Value = -176831.496
$ws.Range("H113").Value = 11743.333
$ws.Range("I113").Value = 11901.818
$ws.Range("K113").Value = 11901.818
$ws.Range("M113").Value = -9731.817999999999
$ws.Range("H136").Value = 13892258
$ws.Range("I136").Value = 7816102.5
$ws.Range("K136").Value = 23448307.5
$ws.Range("M136").Value = -23445757.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 285
$ws.Range("I113").Value = 313.33334
$ws.Range("J113").Value = 200
$ws.Range("K113").Value = 940.0000200000001
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1229.99998
$ws.Range("N113").Value = -4940
$ws.Range("H136").Value = 3118518
$ws.Range("I136").Value = 1740584.5
$ws.Range("K136").Value = 5221753.5
$ws.Range("M136").Value = -5219203.5
